$d = $word.ActiveDocument

# --- Table 1: Sven Johansson reference block ---
# Split the single "Phone:" row (which currently concatenates two numbers)
# into "Phone 1:" / "Phone 2:" rows.
$t1 = $d.Tables.Item(1)
$phoneRow = $t1.Rows.Item($t1.Rows.Count)
$phoneRow.Cells.Item(1).Range.Text = "Phone 1:"
$phoneRow.Cells.Item(2).Range.Text = "+46455 38 57 10"

$newRow1 = $t1.Rows.Add()
$newRow1.Cells.Item(1).Range.Text = "Phone 2:"
$newRow1.Cells.Item(2).Range.Text = "+4670 887 87 08"

# --- Table 2: Markus Fiedler reference block ---
# Append two new rows with phone numbers after the existing "Email:" row.
$t2 = $d.Tables.Item(2)

$newRow2 = $t2.Rows.Add()
$newRow2.Cells.Item(1).Range.Text = "Phone 1:"
$newRow2.Cells.Item(2).Range.Text = "+46455 38 59 23"

$newRow3 = $t2.Rows.Add()
$newRow3.Cells.Item(1).Range.Text = "Phone 2:"
$newRow3.Cells.Item(2).Range.Text = "+4670 853 73 39"
